$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the percentage-like return figures (columns B:E, data rows) are
# stored as literal text rather than being auto-converted to numbers.
$ws.Range("B3:E40").NumberFormat = "@"

# Row 1 = new row index, in order (row 1 .. row 40 of the final layout)
$data = @(
  @('INDEX RETURNS', 'Unnamed: 1', 'Unnamed: 2', 'Unnamed: 3', 'Unnamed: 4'),
  @($null, '1 WEEK', 'MTD', 'QTD', 'YTD'),
  @('EQUITIES', $null, $null, $null, $null),
  @('S&P 500', '3.31%', '3.71%', '13.27%', '-9.34%'),
  @('DJ Industrial Average', '2.99%', '2.87%', '9.89%', '-5.98%'),
  @('Russell 2000', '4.97%', '7.03%', '18.21%', '-9.48%'),
  @('Russell Midcap', '4.19%', '4.91%', '15.26%', '-9.60%'),
  @('STOXX Europe 50 (€)', '1.42%', '1.92%', '9.53%', '-9.58%'),
  @('STOXX Europe 600 (€)†', '1.29%', '0.74%', '8.54%', '-7.71%'),
  @('MSCI EAFE Small Cap', '2.99%', '2.47%', '9.25%', '-17.47%'),
  @('FTSE 100 (£)', '1.18%', '1.54%', '5.27%', '4.25%'),
  @('DAX (€)', '1.63%', '2.31%', '7.92%', '-13.15%'),
  @('FTSE MIB (€)', '1.70%', '2.52%', '8.35%', '-12.89%'),
  @('CAC 40 (€)†', '1.26%', '1.63%', '10.76%', '-6.45%'),
  @('SWISS MKT (CHF)', '0.05%', '-0.16%', '3.60%', '-11.22%'),
  @('TOPIX (¥)', '1.34%', '1.69%', '5.48%', '0.44%'),
  @('Hang Seng (HKD)', '-0.11%', '0.13%', '-7.20%', '-11.67%'),
  @('MSCI World', '3.06%', '3.30%', '11.54%', '-11.09%'),
  @('MSCI China Free†', '-0.38%', '-0.01%', '-8.93%', '-15.12%'),
  @('MSCI EAFE', '2.17%', '1.51%', '6.58%', '-13.95%'),
  @('MSCI EM', '1.66%', '2.65%', '2.49%', '-15.42%'),
  @('MSCI Brazil (BRL)', '6.16%', '9.11%', '14.45%', '10.53%'),
  @('MSCI India (INR)', '1.69%', '3.48%', '13.60%', '2.48%'),
  @('MSCI Russia (RUB)', '0.00%', '0.00%', '0.00%', '-100.00%'),
  @('FIXED INCOME', $null, $null, $null, $null),
  @('Bloomberg Aggregate', '0.24%', '-0.80%', '1.63%', '-8.89%'),
  @('Bloomberg Euro Aggregate', '0.32%', '-0.51%', '1.05%', '-18.37%'),
  @('Bloomberg US High Yield', '0.94%', '1.60%', '7.60%', '-7.67%'),
  @('Bloomberg Euro High Yield (€)', '0.75%', '1.66%', '6.83%', '-8.59%'),
  @('Bloomberg Muni Aggregate', '-0.08%', '-0.19%', '2.44%', '-6.76%'),
  @('Bloomberg TIPS', '0.18%', '-1.41%', '3.06%', '-6.90%'),
  @('JPM EMBI Glbl. Divers.', '1.51%', '2.43%', '5.39%', '-16.01%'),
  @('JPM GBI-EM Glbl. Divers.', '2.34%', '3.09%', '3.40%', '-11.63%'),
  @('OTHER', $null, $null, $null, $null),
  @('DJ US Real Estate', '4.28%', '2.66%', '11.77%', '-10.61%'),
  @('FTSE EPRA/NAREIT Dvlpd. Ex-US', '2.89%', '0.83%', '7.21%', '-14.36%'),
  @('S&P GSCI', '4.45%', '-1.83%', '-1.87%', '33.26%'),
  @('Alerian MLP *', '4.95%', '2.14%', '14.90%', '26.43%'),
  @('US Dollar Index', '-0.91%', '-0.25%', '0.90%', '10.59%'),
  @('VIX', '-7.66%', '-8.44%', '-31.97%', '13.41%')
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 1
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $val = $row[$c]
    if ($val -ne $null) {
      $ws.Cells.Item($r, $c + 1).Value = $val
    } else {
      $ws.Cells.Item($r, $c + 1).ClearContents()
    }
  }
}

# Remove now-unused trailing columns F and G (old layout had 7 columns, new has 5)
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(6).Delete()

Write-Host "Applied index_returns restructuring"
